$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U89"), $null, 1)
Write-Host ("Table name: " + $lo.Name)
Write-Host ("HeaderRowRange: " + $lo.HeaderRowRange.Address())
$lo.HeaderRowRange.Font.Bold = $false
Write-Host "done"
